{"js": "// Fix citation years/details in the \"powe(R)OC Testing Results\" document.\n//\n// The author corrected two references to the Colloff et al. (2021) paper\n// to \"2021a\" (to disambiguate it from the \"2021b\" Colloff et al. paper\n// referenced later in the document), and added the specific sample sizes\n// that were simulated.\n//\n// 1. \"...condition from Colloff et al., 2021 as a base)...\"\n//        -> \"...condition from Colloff et al., 2021a as a base)...\"\n// 2. \"...3 sample sizes, again using...\"\n//        -> \"...3 sample sizes (1000, 3000, 5000), again using...\"\n// 3. \"...condition data from Colloff et al. (2021) as a base...\"\n//        -> \"...condition data from Colloff et al. (2021a) as a base...\"\n\nasync function replaceOnce(searchText, replacementText) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${searchText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Target the smallest unambiguous text spans so existing formatting/run\n// boundaries outside the edit are left untouched, and non-breaking spaces\n// elsewhere in the surrounding text are preserved.\nawait replaceOnce(\"2021 as a base\", \"2021a as a base\");\nawait replaceOnce(\"3 sample sizes, again\", \"3 sample sizes (1000, 3000, 5000), again\");\nawait replaceOnce(\"(2021) as a base\", \"(2021a) as a base\");\n", "ps1": "# Fix citation years/details in the \"powe(R)OC Testing Results\" document.\n#\n# The author corrected two references to the Colloff et al. (2021) paper\n# to \"2021a\" (to disambiguate it from the \"2021b\" Colloff et al. paper\n# referenced later in the document), and added the specific sample sizes\n# that were simulated.\n#\n# 1. \"...condition from Colloff et al., 2021 as a base)...\"\n#        -> \"...condition from Colloff et al., 2021a as a base)...\"\n# 2. \"...3 sample sizes, again using...\"\n#        -> \"...3 sample sizes (1000, 3000, 5000), again using...\"\n# 3. \"...condition data from Colloff et al. (2021) as a base...\"\n#        -> \"...condition data from Colloff et al. (2021a) as a base...\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once([string]$searchText, [string]$replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $found = $find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Could not find text: $searchText\"\n    }\n}\n\n# Target the smallest unambiguous text spans so non-breaking spaces and\n# other text elsewhere in the document are left untouched.\nReplace-Once \"2021 as a base\" \"2021a as a base\"\nReplace-Once \"3 sample sizes, again\" \"3 sample sizes (1000, 3000, 5000), again\"\nReplace-Once \"(2021) as a base\" \"(2021a) as a base\"\n"}
